$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.232.79'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '1.655.97'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').Value = "'219.48"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = "'0.5231"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('D8').Value = "'0.2663"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = "'0.06360"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('D10').Value = "'20.61"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.54%  '
$ws.Range('D11').Value = "'0.07722"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = '1.653.95'
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').Value = '1.884.32'
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('D17').Value = "'65.42"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').Value = '26.225.38'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = "'4.686"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.68%  '
$ws.Range('D21').Value = "'191.87"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.43%  '
$ws.Range('D22').Value = "'10.37"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('D23').Value = "'6.011"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.02%  '
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').Value = "'143.42"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.20%  '
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('D27').Value = "'7.276"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').Value = "'15.95"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.23%  '
$ws.Range('D29').Value = "'1.502"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('E30').Value = '  -4.65%  '
$ws.Range('D31').Value = "'1.276"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('E32').Value = '  -1.64%  '
$ws.Range('D33').Value = "'3.356"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').Value = "'1.579"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('D35').Value = "'2.803"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.415"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.81%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = "'0.9455"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('D38').Value = "'0.5746"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('D39').Value = "'0.01601"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('D40').Value = "'5.902"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').Value = "'2.570"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').Value = "'0.8460"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D44').Value = '1.020.18'
$ws.Range('E44').Value = '  -5.54%  '
$ws.Range('D45').Value = "'101.27"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('D46').Value = '1.795.09'
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').Value = "'58.35"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D50').Value = "'0.05318"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.98%  '
$ws.Range('D51').Value = "'0.4351"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.33%  '
